$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host ("Sheet name: " + $ws.Name)
Write-Host ("A1 value: " + $ws.Range("A1").Value)
Write-Host ("A11 value: " + $ws.Range("A11").Value)
Write-Host ("A11 style: " + $ws.Range("A11").Style)
$ws.Rows("11:13").Insert()
Write-Host "After insert"
Write-Host ("A11 value: " + $ws.Range("A11").Value)
Write-Host ("A14 value: " + $ws.Range("A14").Value)
